$wb = $excel.ActiveWorkbook

# ---- Sheet "Metadata" ----
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refreshed publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty; now populated
$meta.Range("B9").Value = "Alvearie Team"

# Remove the duplicated "Contact" / "No display for ContactDetail" row (old row 11)
$meta.Rows.Item(11).Delete()

# Old row 10 ("Contact" / "No display for ContactDetail") becomes the new
# "Jurisdiction" / "United States of America" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# ---- Sheet "Elements" ----
$elements = $wb.Worksheets.Item("Elements")

# Root "Extension" element row (row 2): Short/Definition customized away from
# the generic placeholder text to describe this specific extension
$elements.Range("K2").Value = "Employee Job Location"
$elements.Range("L2").Value = "Code for the physical location where the employee works"
